$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.107.11"
$ws.Range("E2").Value = "  +0.37%  "
$ws.Range("D3").Value = "2.482.29"
$ws.Range("E3").Value = "  +0.28%  "
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").Value = "'585.55"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.24%  "
$ws.Range("D6").Value = "'171.85"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +2.31%  "
$ws.Range("E7").Value = "  -0.08%  "
$ws.Range("E8").Value = "  -0.13%  "
$ws.Range("D9").Value = "2.481.24"
$ws.Range("E9").Value = "  +0.26%  "
$ws.Range("E10").Value = "  +1.50%  "
$ws.Range("E11").Value = "  +0.02%  "
$ws.Range("E12").Value = "  -0.11%  "
$ws.Range("D13").Value = "'0.330"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -1.40%  "
$ws.Range("D14").Value = "2.934.26"
$ws.Range("E14").Value = "  -0.94%  "
$ws.Range("D15").Value = "'25.42"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -1.54%  "
$ws.Range("D16").Value = "67.013.56"
$ws.Range("E16").Value = "  +0.46%  "
$ws.Range("E17").Value = "  -1.25%  "
$ws.Range("D18").Value = "2.477.20"
$ws.Range("E18").Value = "  +0.52%  "
$ws.Range("E19").Value = "  -4.24%  "
$ws.Range("D20").Value = "'7.40"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -4.43%  "
$ws.Range("D21").Value = "'350.30"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -3.29%  "
$ws.Range("E22").Value = "  -0.97%  "
$ws.Range("E23").Value = "  -0.11%  "
$ws.Range("E24").Value = "  -3.35%  "
$ws.Range("E25").Value = "  -4.53%  "
$ws.Range("E26").Value = "  -2.19%  "
$ws.Range("D27").Value = "'9.26"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.55%  "
$ws.Range("D28").Value = "'1.00"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.09%  "
$ws.Range("E30").Value = "  -2.50%  "
$ws.Range("D31").Value = "'508.81"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -1.19%  "
$ws.Range("D32").Value = "'7.70"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -5.16%  "
$ws.Range("E33").Value = "  -2.95%  "
$ws.Range("E34").Value = "  -3.58%  "
$ws.Range("D35").Value = "'0.999"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.12%  "
$ws.Range("D36").Value = "'159.29"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +1.91%  "
$ws.Range("D37").Value = "'0.117"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -6.72%  "
$ws.Range("E38").Value = "  +0.77%  "
$ws.Range("E39").Value = "  -3.61%  "
$ws.Range("E41").Value = "  -0.14%  "
$ws.Range("E42").Value = "  -3.32%  "
$ws.Range("E43").Value = "  -2.48%  "
$ws.Range("E44").Value = "  -1.19%  "
$ws.Range("E45").Value = "  -2.30%  "
$ws.Range("E46").Value = "  -0.89%  "
$ws.Range("D47").Value = "'142.94"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.27%  "
$ws.Range("E48").Value = "  -4.15%  "
$ws.Range("E49").Value = "  -4.14%  "
$ws.Range("E50").Value = "  -5.98%  "
$ws.Range("E51").Value = "  -1.21%  "

Write-Output "Applied changes to cells"